$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update part number text in B7 from "74HC138PW" to "74HC138PWR"
$ws.Range("B7").Value = "74HC138PWR"

# Update the active selection to B7 (as recorded in the saved workbook view)
$ws.Range("B7").Select()
